# Update cryptos.xlsx price/volume snapshot data.
#
# Each "Price" (column D) and "Volume(1h)" (column E) cell in the source
# workbook is stored as TEXT (the sheet was produced by a scraper that
# writes things like "64.145.51" or "  -3.19%  " as literal strings, not
# numbers). When a value being written still parses as a plain number
# (e.g. "1.00", "608.85", "0.0000259"), Excel's COM automation would
# normally coerce the assignment to a numeric cell and silently drop the
# formatting (trailing zeros, leading zeros, etc.). To keep those cells
# as text - matching the original authoring - we flip the cell to the
# Text number format before writing the value, then restore the cell
# style to Normal so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.099.11'
$ws.Range("E2").Value = '  -3.20%  '
$ws.Range("D3").Value = '3.130.19'
$ws.Range("E3").Value = '  -2.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.55%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '3.126.16'
$ws.Range("E8").Value = '  -2.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.70%  '
$ws.Range("E10").Value = '  -5.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.55'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.477'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.56%  '
$ws.Range("D15").Value = '3.656.96'
$ws.Range("E15").Value = '  -1.95%  '
$ws.Range("D16").Value = '64.204.71'
$ws.Range("E16").Value = '  -3.26%  '
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").Value = '3.137.95'
$ws.Range("E18").Value = '  -2.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.710'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.125'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.71'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.73'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.48%  '
$ws.Range("E35").Value = '  -5.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.09'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.11%  '
$ws.Range("D39").Value = '0.0₃0749'
$ws.Range("E39").Value = '  -3.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '450.07'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.124'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0400'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("D44").Value = '2.874.28'
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.272'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.55%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.46%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.115'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.51%  '
